$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the scrape timestamp for every existing data row (2..410) to the new
# crawl time. Row 411 still exists at this point and is handled (removed)
# further below, so this only needs to cover 2..410 per the diff.
$ws.Range("O2:O410").Value2 = "2022-12-13 20:49:22"

# Product-availability text tweaks captured by the diff.
$ws.Range("M79").Value2 = "Betty Bossi Butterblätterteig eckig ausgewallt 25x42cm - Online kein Bestand 3.70 Schweizer Franken"
$ws.Range("M114").Value2 = "Mini Panettone 20% pro 3 Aktion 1.60 Schweizer Franken"

# The last row (411, Roland Knäckebrot Zwerghirse) was dropped from this crawl
# snapshot entirely.
$ws.Rows.Item(411).Delete()
